$wb = $excel.ActiveWorkbook

# Sheets (by name, per workbook.xml):
#   "limit"                  -> sheet1
#   "offset"                 -> sheet2
#   "inTableOptions"         -> sheet3
#   "optionOnColumn"         -> sheet4
#   "optionForTableAndColumn"-> sheet5

# Rename converter option "type"/"string" to "readAs"/"text" across the workbook.

$wsLimit = $wb.Worksheets.Item("limit")
$wsLimit.Range("C1").Value = "options?limit=5#aaa?readAs=text"

$wsOffset = $wb.Worksheets.Item("offset")
$wsOffset.Range("C1").Value = "optionsOffset?limit=5&offset=2#aaa?readAs=text"

$wsInTable = $wb.Worksheets.Item("inTableOptions")
$wsInTable.Range("C3").Value = "text"
$wsInTable.Range("A3").Value = "single?readAs"
$wsInTable.Range("E8").Value = "text"
$wsInTable.Range("A8").Value = "optionsInTable?readAs"

# Update selections / active sheet to match the saved workbook state.
$wsInTable.Select()
$wsInTable.Range("A9").Select()
$wsInTable.Activate()
